$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.782258064516129
$ws.Range("C2").Value = 0.717305151915456
$ws.Range("D2").Value = 0.801104972375691
$ws.Range("E2").Value = 0.661490683229814
$ws.Range("F2").Value = 0.541808550889141
